# UC3.3_TC1.xlsx - update evaluations (test code generation module)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")
$ws.Activate()

# Assertion validity row: the test now fails, with a different note
$ws.Range("B7").Value = "no"
$ws.Range("C7").Value = "Functionality not working on portal"

# Updated Code BLEU score + the breakdown note that embeds it
$ws.Range("B12").Value = 0.3110227473524679
$ws.Range("C12").Value = "{'codebleu': 0.3110227473524679, 'ngram_match_score': 0.1401184800643443, 'weighted_ngram_match_score': 0.15509175170454562, 'syntax_match_score': 0.6336633663366337, 'dataflow_match_score': 0.31521739130434784}"

# Move the active selection on the sheet as recorded by the author
$ws.Range("B8").Select()
